$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 585; existing rows 585..603 shift down to 587..605.
$ws.Rows.Item(585).Insert()
$ws.Rows.Item(585).Insert()

# New row 585
$ws.Cells.Item(585, 1).Value = 7
$ws.Cells.Item(585, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(585, 3).Value = "Ñuble"
$ws.Cells.Item(585, 4).Value = 44939
$ws.Cells.Item(585, 5).Value = 16
$ws.Cells.Item(585, 6).Value = 100112020
$ws.Cells.Item(585, 7).Value = "Tomate"
$ws.Cells.Item(585, 8).Value = "Larga vida"
$ws.Cells.Item(585, 9).Value = "Primera"
$ws.Cells.Item(585, 10).Value = 600
$ws.Cells.Item(585, 11).Value = 13000
$ws.Cells.Item(585, 12).Value = 14000
$ws.Cells.Item(585, 13).Value = 13500
$ws.Cells.Item(585, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(585, 15).Value = "Región del Maule"
$ws.Cells.Item(585, 16).Value = 750
$ws.Cells.Item(585, 17).Value = 18
$ws.Cells.Item(585, 18).Value = "Hortaliza"

# New row 586
$ws.Cells.Item(586, 1).Value = 7
$ws.Cells.Item(586, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(586, 3).Value = "Ñuble"
$ws.Cells.Item(586, 4).Value = 44939
$ws.Cells.Item(586, 5).Value = 16
$ws.Cells.Item(586, 6).Value = 100112020
$ws.Cells.Item(586, 7).Value = "Tomate"
$ws.Cells.Item(586, 8).Value = "Larga vida"
$ws.Cells.Item(586, 9).Value = "Segunda"
$ws.Cells.Item(586, 10).Value = 300
$ws.Cells.Item(586, 11).Value = 11000
$ws.Cells.Item(586, 12).Value = 11000
$ws.Cells.Item(586, 13).Value = 11000
$ws.Cells.Item(586, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(586, 15).Value = "Región del Maule"
$ws.Cells.Item(586, 16).Value = 611
$ws.Cells.Item(586, 17).Value = 18
$ws.Cells.Item(586, 18).Value = "Hortaliza"

# Column D (Fecha) needs the same date/time number format as the rest of the column.
$ws.Cells.Item(585, 4).NumberFormat = $ws.Cells.Item(584, 4).NumberFormat
$ws.Cells.Item(586, 4).NumberFormat = $ws.Cells.Item(584, 4).NumberFormat
